# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (every sheet except the first five
# reference/meta sheets), prepend the worksheet's own name + a space to
# each Step/label value held in column A (rows 2..last used row). Column A
# row 1 is the "Name" header and is left untouched.

$wb = $excel.ActiveWorkbook

# These sheets hold reference data (not protocol steps) and must not be touched.
$skipSheets = @("ZansiJourney", "NRWaves", "PersonalZansi", "PositiveSpin", "ReEngagement")

foreach ($ws in $wb.Worksheets) {
    $sheetName = $ws.Name

    if ($skipSheets -contains $sheetName) {
        continue
    }

    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Formula

        if ($current -eq $null -or $current -eq "") {
            continue
        }

        $prefix = $sheetName + " "
        if ($current.StartsWith($prefix)) {
            continue
        }

        $cell.Formula = $prefix + $current
    }
}
